# Generate Report for Handoff
# Updates the localization-status workbook to reflect a fresh handoff
# generated for "b.md" (source file), for both the zh-cn and de-de locales.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview": columns A=File Name, B=zh-cn, C=de-de
# Row 3 corresponds to b.md; its status moves from
# "Handed back: in sync with en-US" to "Ready for handoff".
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# ---------------------------------------------------------------------
# Sheet "zh-cn": columns A=Source File Name, B=Status, C=Latest Handoff
# File, D=Latest Handoff Datetime, E=Latest Target File, F=Latest
# Handback File, G=Latest Handback DateTime, H=Handoff Reason,
# I=Dependency From.
# Row 3 corresponds to b.md and gets a new handoff file + datetime.
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = "Ready for handoff"
$wsZhCn.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("D3").Value = "2016-03-03 15:03:11"

foreach ($hl in $wsZhCn.Hyperlinks) {
    if ($hl.Range.Address() -eq '$C$3') {
        $hl.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
    }
}

# ---------------------------------------------------------------------
# Sheet "de-de": same column layout as zh-cn.
# Row 3 corresponds to b.md and gets a new handoff file + datetime.
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = "Ready for handoff"
$wsDeDe.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("D3").Value = "2016-03-03 15:03:25"

foreach ($hl in $wsDeDe.Hyperlinks) {
    if ($hl.Range.Address() -eq '$C$3') {
        $hl.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
    }
}
